# Update Slamf7-Slamf7 sheet with new TPM-derived values.
# The underlying data pipeline was re-run with new TPM values, which
# collapsed the 4-row cross product (ECs/MuSCs x ECs/MuSCs) down to a
# single MuSCs-MuSCs self-pair row, and removed the now-unused "ECs"
# string from the shared string table (handled automatically by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3-5 (old rows for ECs-MuSCs, MuSCs-ECs, MuSCs-MuSCs-old)
# leaving only the header row and the single updated data row.
$ws.Range("A3:T5").EntireRow.Delete()

# Update the remaining data row (row 2) with the new computed values.
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Slamf7"
$ws.Range("C2").Value = "Slamf7"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03958366666666666
$ws.Range("H2").Value = 0.118751
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03958366666666666
$ws.Range("N2").Value = 0.118751
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.001566866666777777
$ws.Range("R2").Value = 0.014101800001
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
